# Update bus voltage magnitude results (vm_pu) for Case_3_23 - "case with 380 kV done"
# Slack bus voltage setpoint changed from 1.05 pu to 1.02 pu, and all dependent
# bus voltage results (columns B-F, I-N, rows 2-25) were recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037008187272881
$ws.Range("D2").Value = 1.045203459796238
$ws.Range("E2").Value = 1.046113121317894
$ws.Range("F2").Value = 1.057626833267839
$ws.Range("I2").Value = 1.041006179064552
$ws.Range("J2").Value = 1.042113398590972
$ws.Range("K2").Value = 1.047972358742998
$ws.Range("L2").Value = 1.04887946797452
$ws.Range("M2").Value = 1.060361286356204
$ws.Range("N2").Value = 1.017985882721226
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037754295821383
$ws.Range("D3").Value = 1.045772430693646
$ws.Range("E3").Value = 1.046763536611334
$ws.Range("F3").Value = 1.058326402752791
$ws.Range("I3").Value = 1.041157168498658
$ws.Range("J3").Value = 1.042504731225998
$ws.Range("K3").Value = 1.048353354407512
$ws.Range("L3").Value = 1.049341879369079
$ws.Range("M3").Value = 1.060875019614482
$ws.Range("N3").Value = 1.018116536212829
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038237914142969
$ws.Range("D4").Value = 1.046141338139337
$ws.Range("E4").Value = 1.047185513464503
$ws.Range("F4").Value = 1.058780225670402
$ws.Range("I4").Value = 1.041254127487823
$ws.Range("J4").Value = 1.042758077173767
$ws.Range("K4").Value = 1.048599924094189
$ws.Range("L4").Value = 1.04964151640736
$ws.Range("M4").Value = 1.0612079144454
$ws.Range("N4").Value = 1.018201095552936
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038441425392033
$ws.Range("D5").Value = 1.046296603210284
$ws.Range("E5").Value = 1.047363177241426
$ws.Range("F5").Value = 1.058971286970401
$ws.Range("I5").Value = 1.041294710780563
$ws.Range("J5").Value = 1.042864612948132
$ws.Range("K5").Value = 1.04870359012121
$ws.Range("L5").Value = 1.049767584170745
$ws.Range("M5").Value = 1.061347975353248
$ws.Range("N5").Value = 1.018236648059422
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038475607382818
$ws.Range("D6").Value = 1.046322683182885
$ws.Range("E6").Value = 1.047393023231464
$ws.Range("F6").Value = 1.059003383012809
$ws.Range("I6").Value = 1.041301514416999
$ws.Range("J6").Value = 1.042882502441319
$ws.Range("K6").Value = 1.04872099653793
$ws.Range("L6").Value = 1.049788757340523
$ws.Range("M6").Value = 1.061371498706411
$ws.Range("N6").Value = 1.018242617687626
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038240632694426
$ws.Range("D7").Value = 1.046143412108316
$ws.Range("E7").Value = 1.04718788637917
$ws.Range("F7").Value = 1.058782777566523
$ws.Range("I7").Value = 1.041254670465235
$ws.Range("J7").Value = 1.042759500597268
$ws.Range("K7").Value = 1.04860130925459
$ws.Range("L7").Value = 1.049643200538707
$ws.Range("M7").Value = 1.061209785508611
$ws.Range("N7").Value = 1.018201570593123
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037260163757764
$ws.Range("D8").Value = 1.045395590651774
$ws.Range("E8").Value = 1.04633270012361
$ws.Range("F8").Value = 1.057863015591586
$ws.Range("I8").Value = 1.041057359651761
$ws.Range("J8").Value = 1.042245624031661
$ws.Range("K8").Value = 1.04810110885539
$ws.Range("L8").Value = 1.049035652902244
$ws.Range("M8").Value = 1.0605348054644
$ws.Range("N8").Value = 1.018030033695462
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03553894868167
$ws.Range("D9").Value = 1.044083630929036
$ws.Range("E9").Value = 1.04483438454398
$ws.Range("F9").Value = 1.056251222127377
$ws.Range("I9").Value = 1.040704031196312
$ws.Range("J9").Value = 1.041341146862657
$ws.Range("K9").Value = 1.047220062408163
$ws.Range("L9").Value = 1.047968410205289
$ws.Range("M9").Value = 1.059349122811271
$ws.Range("N9").Value = 1.017727922162002
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034395964617496
$ws.Range("D10").Value = 1.043213010340663
$ws.Range("E10").Value = 1.043841439609811
$ws.Range("F10").Value = 1.055182844130237
$ws.Range("I10").Value = 1.040464734140937
$ws.Range("J10").Value = 1.040738944624008
$ws.Range("K10").Value = 1.046633034094863
$ws.Range("L10").Value = 1.047259253619486
$ws.Range("M10").Value = 1.058561276207512
$ws.Range("N10").Value = 1.017526650915423
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033902131994614
$ws.Range("D11").Value = 1.042837001156774
$ws.Range("E11").Value = 1.043412917850508
$ws.Range("F11").Value = 1.054721712755062
$ws.Range("I11").Value = 1.040360237934317
$ws.Range("J11").Value = 1.04047838785436
$ws.Range("K11").Value = 1.046378942498025
$ws.Range("L11").Value = 1.046952756287741
$ws.Range("M11").Value = 1.058220771685222
$ws.Range("N11").Value = 1.017439536965566
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033718866042916
$ws.Range("D12").Value = 1.042697483243766
$ws.Range("E12").Value = 1.043253962865407
$ws.Range("F12").Value = 1.054550653277812
$ws.Range("I12").Value = 1.040321292155192
$ws.Range("J12").Value = 1.040381637085501
$ws.Range("K12").Value = 1.04628457742136
$ws.Range("L12").Value = 1.046838997175922
$ws.Range("M12").Value = 1.058094390919105
$ws.Range("N12").Value = 1.017407185166802
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033758169689011
$ws.Range("D13").Value = 1.042727403549539
$ws.Range("E13").Value = 1.043288049390187
$ws.Range("F13").Value = 1.054587335885097
$ws.Range("I13").Value = 1.040329652082439
$ws.Range("J13").Value = 1.040402388999822
$ws.Range("K13").Value = 1.046304818309022
$ws.Range("L13").Value = 1.046863394898474
$ws.Range("M13").Value = 1.058121495567737
$ws.Range("N13").Value = 1.017414124448549
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033886979785798
$ws.Range("D14").Value = 1.042825465518141
$ws.Range("E14").Value = 1.043399774140613
$ws.Range("F14").Value = 1.054707568316521
$ws.Range("I14").Value = 1.040357021336451
$ws.Range("J14").Value = 1.040470389752258
$ws.Range("K14").Value = 1.046371141919214
$ws.Range("L14").Value = 1.046943351126426
$ws.Range("M14").Value = 1.058210323006219
$ws.Range("N14").Value = 1.0174368626263
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033966365929007
$ws.Range("D15").Value = 1.042885904492844
$ws.Range("E15").Value = 1.043468640278432
$ws.Range("F15").Value = 1.054781677404596
$ws.Range("I15").Value = 1.040373867071891
$ws.Range("J15").Value = 1.040512291493086
$ws.Range("K15").Value = 1.046412008230371
$ws.Range("L15").Value = 1.046992626461387
$ws.Range("M15").Value = 1.058265065535094
$ws.Range("N15").Value = 1.017450873205225
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034428761748888
$ws.Range("D16").Value = 1.043237985560479
$ws.Range("E16").Value = 1.043869909506596
$ws.Range("F16").Value = 1.055213479357551
$ws.Range("I16").Value = 1.040471650760836
$ws.Range("J16").Value = 1.040756241272794
$ws.Range("K16").Value = 1.046649899458461
$ws.Range("L16").Value = 1.047279607036395
$ws.Range("M16").Value = 1.058583887979412
$ws.Range("N16").Value = 1.017532433222579
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034719102980913
$ws.Range("D17").Value = 1.043459099316366
$ws.Range("E17").Value = 1.044121999514667
$ws.Range("F17").Value = 1.05548473613828
$ws.Range("I17").Value = 1.040532753133542
$ws.Range("J17").Value = 1.040909319303331
$ws.Range("K17").Value = 1.04679914897997
$ws.Range("L17").Value = 1.047459776740192
$ws.Range("M17").Value = 1.058784049085495
$ws.Range("N17").Value = 1.017583604160081
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034888558846096
$ws.Range("D18").Value = 1.043588165187617
$ws.Range("E18").Value = 1.044269177170619
$ws.Range("F18").Value = 1.055643098590611
$ws.Range("I18").Value = 1.040568308247395
$ws.Range("J18").Value = 1.04099862640307
$ws.Range("K18").Value = 1.046886212817997
$ws.Range("L18").Value = 1.047564921771474
$ws.Range("M18").Value = 1.058900861128508
$ws.Range("N18").Value = 1.017613454917374
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034946356616454
$ws.Range("D19").Value = 1.043632189186894
$ws.Range("E19").Value = 1.044319384256135
$ws.Range("F19").Value = 1.055697120270333
$ws.Range("I19").Value = 1.040580417202502
$ws.Range("J19").Value = 1.041029081041233
$ws.Range("K19").Value = 1.046915900837534
$ws.Range("L19").Value = 1.047600782803086
$ws.Range("M19").Value = 1.058940701383186
$ws.Range("N19").Value = 1.017623633851633
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034687941265505
$ws.Range("D20").Value = 1.043435366169152
$ws.Range("E20").Value = 1.044094938371099
$ws.Range("F20").Value = 1.05545561804845
$ws.Range("I20").Value = 1.040526206199406
$ws.Range("J20").Value = 1.04089289347526
$ws.Range("K20").Value = 1.046783134964929
$ws.Range("L20").Value = 1.047440440525589
$ws.Range("M20").Value = 1.058762567322829
$ws.Range("N20").Value = 1.017578113625097
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033849043837606
$ws.Range("D21").Value = 1.04279658459237
$ws.Range("E21").Value = 1.043366867975927
$ws.Range("F21").Value = 1.054672156625071
$ws.Range("I21").Value = 1.040348965389412
$ws.Range("J21").Value = 1.040450364341106
$ws.Range("K21").Value = 1.046351610819329
$ws.Range("L21").Value = 1.046919803577229
$ws.Range("M21").Value = 1.058184162825387
$ws.Range("N21").Value = 1.017430166624356
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033322554170403
$ws.Range("D22").Value = 1.042395818288857
$ws.Range("E22").Value = 1.042910358264061
$ws.Range("F22").Value = 1.054180867724724
$ws.Range("I22").Value = 1.040236768179628
$ws.Range("J22").Value = 1.040172312352786
$ws.Range("K22").Value = 1.046080386941103
$ws.Range("L22").Value = 1.046592966376857
$ws.Range("M22").Value = 1.057821063583884
$ws.Range("N22").Value = 1.017337182605446
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03360156473394
$ws.Range("D23").Value = 1.042608189768606
$ws.Range("E23").Value = 1.043152242799142
$ws.Range("F23").Value = 1.054441184827775
$ws.Range("I23").Value = 1.040296317700512
$ws.Range("J23").Value = 1.04031969508089
$ws.Range("K23").Value = 1.046224158583505
$ws.Range("L23").Value = 1.046766180254543
$ws.Range("M23").Value = 1.058013495024388
$ws.Range("N23").Value = 1.017386471595611
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034702021580034
$ws.Range("D24").Value = 1.043446089866775
$ws.Range("E24").Value = 1.044107165711518
$ws.Range("F24").Value = 1.055468774818476
$ws.Range("I24").Value = 1.040529164739469
$ws.Range("J24").Value = 1.040900315540569
$ws.Range("K24").Value = 1.046790370981046
$ws.Range("L24").Value = 1.04744917755866
$ws.Range("M24").Value = 1.058772273817077
$ws.Range("N24").Value = 1.01758059455042
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035983141316066
$ws.Range("D25").Value = 1.044422104428382
$ws.Range("E25").Value = 1.045220698855064
$ws.Range("F25").Value = 1.056666835244105
$ws.Range("I25").Value = 1.040796038653067
$ws.Range("J25").Value = 1.041574844106985
$ws.Range("K25").Value = 1.047447781166116
$ws.Range("L25").Value = 1.048243912936276
$ws.Range("M25").Value = 1.059655198362541
$ws.Range("N25").Value = 1.017806003375775
